# Fixing cut off issues and stage 6
# 1. Rename the "Stage# Time" / "Time Spent in ... OUT" headers in column B
# 2. Add a new column H "Running Total Lifebase Duration_seconds" with the
#    running (cumulative) total of column G (Lifebase Duration_seconds)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename column B labels (drop "OUT" suffix / add space in "Stage N") ---
$ws.Range("B2").Value  = "Stage 1"
$ws.Range("B3").Value  = "Time Spent in Valgrisenche"
$ws.Range("B4").Value  = "Stage 2"
$ws.Range("B5").Value  = "Time Spent in Cogne"
$ws.Range("B6").Value  = "Stage 3"
$ws.Range("B7").Value  = "Time Spent in Donnas"
$ws.Range("B8").Value  = "Stage 4"
$ws.Range("B9").Value  = "Time Spent in Gressoney"
$ws.Range("B10").Value = "Stage 5"
$ws.Range("B11").Value = "Time Spent in Valtournenche"
$ws.Range("B12").Value = "Stage 6"
$ws.Range("B13").Value = "Time Spent in Ollomont"
$ws.Range("B14").Value = "Stage 7"

# --- Add new column H: running total of column G ---
$ws.Range("H1").Value = "Running Total Lifebase Duration_seconds"
# Match the header formatting used by the other header cells (e.g. G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

$runningTotal = 0
for ($row = 2; $row -le 14; $row++) {
    $runningTotal = $runningTotal + $ws.Cells.Item($row, 7).Value()
    $ws.Cells.Item($row, 8).Value = $runningTotal
}
